$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the whole of column B (the "Location" column being discontinued)
# the way a user would before deleting it -- this is what leaves the
# B1:B1048576 column selection behind in the saved view state.
$ws.Columns("B:B").Select() | Out-Null

# Remove column B entirely (values, the shared strings it referenced, and
# its column-width definition) and shift remaining columns left.
$ws.Columns("B:B").Delete()

# Scroll the view down so row 7 is the first visible row (matches the
# author's saved scroll position).
$excel.ActiveWindow.ScrollRow = 7
